$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.724909
$ws.Range("H2").Value = 1.449818
$ws.Range("I2").Value = 0.5784209917375155
$ws.Range("J2").Value = 0.477722171991027
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.684528
$ws.Range("N2").Value = 3.369056
$ws.Range("O2").Value = 0.1245924002783664
$ws.Range("P2").Value = 0.09862679185493829
$ws.Range("Q2").Value = 1.221129507952
$ws.Range("R2").Value = 4.884518031808001
$ws.Range("S2").Value = 0.07206685973197018
$ws.Range("T2").Value = 0.04711620522144805

$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.724909
$ws.Range("H3").Value = 1.449818
$ws.Range("I3").Value = 0.5784209917375155
$ws.Range("J3").Value = 0.477722171991027
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.007174
$ws.Range("N3").Value = 6.021522
$ws.Range("O3").Value = 0.1484562004527854
$ws.Range("P3").Value = 0.1762759054595506
$ws.Range("Q3").Value = 1.455018497166
$ws.Range("R3").Value = 8.730110982996001
$ws.Range("S3").Value = 0.08587018269548355
$ws.Range("T3").Value = 0.08421090842582143

$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.724909
$ws.Range("H4").Value = 1.449818
$ws.Range("I4").Value = 0.5784209917375155
$ws.Range("J4").Value = 0.477722171991027
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.067628
$ws.Range("N4").Value = 6.202884000000001
$ws.Range("O4").Value = 0.1529275473027211
$ws.Range("P4").Value = 0.1815851529830098
$ws.Range("Q4").Value = 1.498842145852
$ws.Range("R4").Value = 8.993052875112001
$ws.Range("S4").Value = 0.08845650357482575
$ws.Range("T4").Value = 0.08674725368436635

$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.724909
$ws.Range("H5").Value = 1.449818
$ws.Range("I5").Value = 0.5784209917375155
$ws.Range("J5").Value = 0.477722171991027
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.716762
$ws.Range("N5").Value = 9.433524
$ws.Range("O5").Value = 0.3488649040691446
$ws.Range("P5").Value = 0.2761599118585636
$ws.Range("Q5").Value = 3.419223224658
$ws.Range("R5").Value = 13.676892898632
$ws.Range("S5").Value = 0.2017907837940878
$ws.Range("T5").Value = 0.1319277129099236

$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.724909
$ws.Range("H6").Value = 1.449818
$ws.Range("I6").Value = 0.5784209917375155
$ws.Range("J6").Value = 0.477722171991027
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.2530873333333333
$ws.Range("N6").Value = 0.759262
$ws.Range("O6").Value = 0.01871904672409779
$ws.Range("P6").Value = 0.02222687163328961
$ws.Range("Q6").Value = 0.1834652857193333
$ws.Range("R6").Value = 1.100791714316
$ws.Range("S6").Value = 0.01082748957053353
$ws.Range("T6").Value = 0.01061826939322086

$ws.Range("D7").Value = "Neutrophils"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.724909
$ws.Range("H7").Value = 1.449818
$ws.Range("I7").Value = 0.5784209917375155
$ws.Range("J7").Value = 0.477722171991027
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.791131666666666
$ws.Range("N7").Value = 8.373394999999999
$ws.Range("O7").Value = 0.2064399011728847
$ws.Range("P7").Value = 0.245125366210648
$ws.Range("Q7").Value = 2.023316465351666
$ws.Range("R7").Value = 12.13989879211
$ws.Range("S7").Value = 0.1194091723706147
$ws.Range("T7").Value = 0.1171018223562466

$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.528346
$ws.Range("H8").Value = 1.585038
$ws.Range("I8").Value = 0.4215790082624845
$ws.Range("J8").Value = 0.522277828008973
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.684528
$ws.Range("N8").Value = 3.369056
$ws.Range("O8").Value = 0.1245924002783664
$ws.Range("P8").Value = 0.09862679185493829
$ws.Range("Q8").Value = 0.890013630688
$ws.Range("R8").Value = 5.340081784127999
$ws.Range("S8").Value = 0.05252554054639618
$ws.Range("T8").Value = 0.05151058663349024

$ws.Range("D9").Value = "ECs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.528346
$ws.Range("H9").Value = 1.585038
$ws.Range("I9").Value = 0.4215790082624845
$ws.Range("J9").Value = 0.522277828008973
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.007174
$ws.Range("N9").Value = 6.021522
$ws.Range("O9").Value = 0.1484562004527854
$ws.Range("P9").Value = 0.1762759054595506
$ws.Range("Q9").Value = 1.060482354204
$ws.Range("R9").Value = 9.544341187836
$ws.Range("S9").Value = 0.06258601775730188
$ws.Range("T9").Value = 0.09206499703372917

$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.528346
$ws.Range("H10").Value = 1.585038
$ws.Range("I10").Value = 0.4215790082624845
$ws.Range("J10").Value = 0.522277828008973
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.067628
$ws.Range("N10").Value = 6.202884000000001
$ws.Range("O10").Value = 0.1529275473027211
$ws.Range("P10").Value = 0.1815851529830098
$ws.Range("Q10").Value = 1.092422983288
$ws.Range("R10").Value = 9.831806849592001
$ws.Range("S10").Value = 0.06447104372789533
$ws.Range("T10").Value = 0.09483789929864346

$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.528346
$ws.Range("H11").Value = 1.585038
$ws.Range("I11").Value = 0.4215790082624845
$ws.Range("J11").Value = 0.522277828008973
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.716762
$ws.Range("N11").Value = 9.433524
$ws.Range("O11").Value = 0.3488649040691446
$ws.Range("P11").Value = 0.2761599118585636
$ws.Range("Q11").Value = 2.492082335652
$ws.Range("R11").Value = 14.952494013912
$ws.Range("S11").Value = 0.1470741202750567
$ws.Range("T11").Value = 0.14423219894864

$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.528346
$ws.Range("H12").Value = 1.585038
$ws.Range("I12").Value = 0.4215790082624845
$ws.Range("J12").Value = 0.522277828008973
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.2530873333333333
$ws.Range("N12").Value = 0.759262
$ws.Range("O12").Value = 0.01871904672409779
$ws.Range("P12").Value = 0.02222687163328961
$ws.Range("Q12").Value = 0.1337176802173333
$ws.Range("R12").Value = 1.203459121956
$ws.Range("S12").Value = 0.007891557153564254
$ws.Range("T12").Value = 0.01160860224006875

$ws.Range("D13").Value = "Neutrophils"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.528346
$ws.Range("H13").Value = 1.585038
$ws.Range("I13").Value = 0.4215790082624845
$ws.Range("J13").Value = 0.522277828008973
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.791131666666666
$ws.Range("N13").Value = 8.373394999999999
$ws.Range("O13").Value = 0.2064399011728847
$ws.Range("P13").Value = 0.245125366210648
$ws.Range("Q13").Value = 1.474683251556666
$ws.Range("R13").Value = 13.27214926401
$ws.Range("S13").Value = 0.08703072880227004
$ws.Range("T13").Value = 0.1280235438544013
